$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: D12 = "Rien" (existing style kept), E12 = "ignore" (default style)
$ws.Range("D12").Value = "Rien"
$ws.Range("E12").Value = "ignore"

# Row 13: D13 = "Rien" (existing style kept), E13 = "ignore" (default style)
$ws.Range("D13").Value = "Rien"
$ws.Range("E13").Value = "ignore"

# Row 30: add E30 = "ignore", formatted like the rest of the row
$ws.Range("C30").Copy($ws.Range("E30"))
$ws.Range("E30").Value = "ignore"

# Row 31: D31 gets the "Rien pour le moment..." placeholder (format copied from C31),
# and E31 = "ignore"
$ws.Range("C31").Copy($ws.Range("D31"))
$ws.Range("D31").Value = "Rien pour le moment..."
$ws.Range("C31").Copy($ws.Range("E31"))
$ws.Range("E31").Value = "ignore"

# Row 36: D36 gets the placeholder text, E36 = "ignore"
$ws.Range("D36").Value = "Rien pour le moment..."
$ws.Range("C36").Copy($ws.Range("E36"))
$ws.Range("E36").Value = "ignore"

# Row 37: D37 gets the placeholder text, E37 = "ignore"
$ws.Range("D37").Value = "Rien pour le moment..."
$ws.Range("C37").Copy($ws.Range("E37"))
$ws.Range("E37").Value = "ignore"

# Row 38: C38 formatting aligns with the rest of the row (copied from C37), D38 gets
# the placeholder text, E38 = "ignore"
$ws.Range("C37").Copy($ws.Range("C38"))
$ws.Range("C38").Value = "Rien pour le moment..."
$ws.Range("D38").Value = "Rien pour le moment..."
$ws.Range("C38").Copy($ws.Range("E38"))
$ws.Range("E38").Value = "ignore"

$excel.CutCopyMode = $false

# Selection moves onto the newly-filled E36:E38 range
$ws.Range("E36:E38").Select()
